$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose displayed "Price" (column D) text looks like a plain number
# need to be forced to Text format first, otherwise Excel's COM layer will
# auto-convert the assigned string into a numeric value (losing the exact
# formatting, e.g. trailing zeros / precision). NumberFormat "@" forces text
# entry, then ClearFormats() drops the now-unneeded style again so the cell
# ends up with the same (default) style index it started with.
$priceUpdates = @{
    'D2'  = '43.960.86'
    'D3'  = '2.253.42'
    'D4'  = '1.01'
    'D5'  = '231.66'
    'D6'  = '0.648'
    'D7'  = '63.80'
    'D10' = '0.0977'
    'D11' = '56.76'
    'D12' = '26.66'
    'D13' = '0.106'
    'D14' = '2.590.86'
    'D15' = '15.55'
    'D16' = '6.11'
    'D17' = '0.833'
    'D18' = '2.268.99'
    'D19' = '43.790.33'
    'D21' = '73.23'
    'D22' = '6.05'
    'D23' = '250.39'
    'D24' = '0.999'
    'D26' = '3.37'
    'D28' = '9.97'
    'D29' = '170.85'
    'D30' = '20.89'
    'D32' = '1.38'
    'D34' = '0.0706'
    'D36' = '4.90'
    'D37' = '3.67'
    'D38' = '6.43'
    'D42' = '0.000219'
    'D43' = '17.36'
    'D45' = '4.43'
    'D46' = '97.49'
    'D47' = '1.20'
    'D50' = '1.436.12'
    'D51' = '2.77'
}

foreach ($cellRef in $priceUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = '@'
    $range.Value = $priceUpdates[$cellRef]
    $range.ClearFormats()
}

# Column E ("Volume(1h)") values are percent strings with surrounding
# whitespace (e.g. "  -0.87%  "); Excel can't coerce these to numbers so a
# plain .Value assignment keeps them as text already.
$volumeUpdates = @{
    'E2'  = '  -0.87%  '
    'E3'  = '  -1.23%  '
    'E4'  = '  +0.54%  '
    'E5'  = '  +0.34%  '
    'E6'  = '  +3.47%  '
    'E7'  = '  +4.84%  '
    'E8'  = '  -0.06%  '
    'E9'  = '  +7.68%  '
    'E10' = '  +3.63%  '
    'E11' = '  -2.26%  '
    'E12' = '  +8.93%  '
    'E13' = '  +2.08%  '
    'E14' = '  -1.13%  '
    'E15' = '  -0.01%  '
    'E16' = '  +5.02%  '
    'E17' = '  +3.17%  '
    'E18' = '  -0.47%  '
    'E19' = '  -0.89%  '
    'E20' = '  +5.60%  '
    'E21' = '  +0.03%  '
    'E22' = '  -2.70%  '
    'E23' = '  -1.39%  '
    'E24' = '  -0.10%  '
    'E25' = '  -4.99%  '
    'E26' = '  +24.74%  '
    'E27' = '  -4.37%  '
    'E28' = '  +1.50%  '
    'E29' = '  -0.31%  '
    'E30' = '  +1.37%  '
    'E31' = '  -1.64%  '
    'E32' = '  -2.43%  '
    'E33' = '  +3.66%  '
    'E34' = '  +7.57%  '
    'E35' = '  +2.15%  '
    'E36' = '  -2.74%  '
    'E37' = '  +2.26%  '
    'E38' = '  -0.84%  '
    'E39' = '  -4.29%  '
    'E40' = '  +3.78%  '
    'E41' = '  -0.03%  '
    'E42' = '  -1.33%  '
    'E43' = '  +4.46%  '
    'E44' = '  -5.76%  '
    'E45' = '  -1.37%  '
    'E46' = '  -0.72%  '
    'E47' = '  -0.91%  '
    'E48' = '  -2.24%  '
    'E49' = '  +4.16%  '
    'E50' = '  -2.78%  '
    'E51' = '  +1.37%  '
}

foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}

# Row 51 changed to a different coin entirely (Celestia -> HuobiToken), so
# the name/link text cells need updating too (price/volume handled above).
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'

Write-Output "done"
